# Apply "added state-space simulation data to spreadsheet" edit.
#
# This reshapes the "Single-span" worksheet:
#  - inserts a "damping ratio" row (row 8) under the bridge/beam block
#  - shifts the old "Vehicle (4Hz) Parameters" block down one row (9 -> 10..14)
#  - adds a mirrored vehicle-parameters block in columns E:G next to the
#    bridge/beam block (rows 1-5)
#  - adds two new sections below (rows 16-20 and 22-26) with EI / total
#    mass / damping-ratio and a raw material-parameters block
#  - centers + merges the three/five section header rows
#  - widens column A to fit the longer labels

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Section 1 (rows 1-8): existing Bridge/beam parameters block, plus a
# mirrored Vehicle (4Hz) Parameters block in E:G and a new damping-ratio
# row.
# ------------------------------------------------------------------

# Row 1 header stays "Bridge/beam parameters" (already there); add the
# second header "Vehicle (4Hz) Parameters" alongside it in E1.
$ws.Range("E1").Value = "Vehicle (4Hz) Parameters"

# E2:G5 - copy of the vehicle parameters, referencing F-column values.
$ws.Range("E2").Value = "mass"
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = "slinch"

$ws.Range("E3").Value = "spring k"
$ws.Range("F3").Value = 63165.5
$ws.Range("F3").NumberFormat = "0.00E+00"
$ws.Range("G3").Value = "lb/in"

$ws.Range("E4").Value = "damping"
$ws.Range("F4").Value = 502.65
$ws.Range("G4").Value = "lb-s/in"

$ws.Range("E5").Value = "percent damping"
$ws.Range("F5").Formula = "=F4/(2*SQRT(F2*F3))"
$ws.Range("F5").NumberFormat = "0.00%"

# New row 8: damping ratio (directly under Total Mass in row 7).
$ws.Range("A8").Value = "damping ratio"
$ws.Range("B8").Value = 0

# ------------------------------------------------------------------
# Section 2 (rows 10-14): the original "Vehicle (4Hz) Parameters" block,
# now shifted one row down (was rows 9-13). Clear the old row 9 header
# cell first since the block moved down into row 10.
# ------------------------------------------------------------------

$ws.Range("A9").ClearContents()

$ws.Range("A10").Value = "Vehicle (4Hz) Parameters"
$ws.Range("B10:C10").ClearContents()

$ws.Range("A11").Value = "mass"
$ws.Range("B11").Style = "Normal"
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = "slinch"

$ws.Range("A12").Value = "spring k"
$ws.Range("B12").Value = 63165.5
$ws.Range("B12").NumberFormat = "0.00E+00"
$ws.Range("C12").Value = "lb/in"

$ws.Range("A13").Value = "damping"
$ws.Range("B13").Style = "Normal"
$ws.Range("B13").Value = 502.65
$ws.Range("C13").Value = "lb-s/in"

$ws.Range("A14").Value = "percent damping"
$ws.Range("B14").Formula = "=B13/(2*SQRT(B11*B12))"
$ws.Range("B14").NumberFormat = "0.00%"

# ------------------------------------------------------------------
# Section 3 (rows 16-20): second Bridge/beam parameters block (EI /
# Total Mass / damping ratio) plus a mirrored Vehicle Parameters block.
# ------------------------------------------------------------------

$ws.Range("A16").Value = "Bridge/beam parameters"
$ws.Range("E16").Value = "Vehicle Parameters"

$ws.Range("A17").Value = "Length"
$ws.Range("B17").Value = 1200
$ws.Range("C17").Value = "in"
$ws.Range("E17").Value = "mass"
$ws.Range("F17").Value = 100
$ws.Range("G17").Value = "slinch"

$ws.Range("A18").Value = "EI"
$ws.Range("B18").Formula = "=B3*B4"
$ws.Range("B18").NumberFormat = "0.00E+00"
$ws.Range("C18").Value = "lb-in^2"
$ws.Range("E18").Value = "spring k"
$ws.Range("F18").Value = 63165.5
$ws.Range("F18").NumberFormat = "0.00E+00"
$ws.Range("G18").Value = "lb/in"

$ws.Range("A19").Value = "Total Mass"
$ws.Range("B19").Formula = "=B7"
$ws.Range("C19").Value = "lb"
$ws.Range("E19").Value = "damping coefficient"
$ws.Range("F19").Value = 502.65
$ws.Range("G19").Value = "lb-s/in"

$ws.Range("A20").Value = "damping ratio"
$ws.Range("B20").Value = 0
$ws.Range("F20").NumberFormat = "0.00%"

# ------------------------------------------------------------------
# Section 4 (rows 22-26): raw material / section parameters block.
# ------------------------------------------------------------------

$ws.Range("A22").Value = "Bridge/beam parameters"

$ws.Range("A23").Value = "Moment of Inertia (I)"
$ws.Range("B23").Value = 1500
$ws.Range("C23").Value = "in^4"

$ws.Range("A24").Value = "Cross sectional area (A)"
$ws.Range("B24").Value = 10
$ws.Range("C24").Value = "in^2"

$ws.Range("A25").Value = "Modulus of Elasticity (E)"
$ws.Range("B25").Value = 5000000000
$ws.Range("B25").NumberFormat = "0.00E+00"
$ws.Range("C25").Value = "psi"

$ws.Range("A26").Value = "Material Density"
$ws.Range("B26").Value = 0.099286
$ws.Range("C26").Value = "slinch/in^3"

# ------------------------------------------------------------------
# Formatting: center + merge the section header rows.
# ------------------------------------------------------------------

$headerRanges = @("A1:C1", "E1:G1", "A16:C16", "E16:G16", "A22:C22")
foreach ($rng in $headerRanges) {
    $ws.Range($rng).HorizontalAlignment = -4108
    $ws.Range($rng).MergeCells = $true
}

# Widen column A to fit the longest label ("Modulus of Elasticity (E)").
$ws.Columns.Item(1).ColumnWidth = 20.7

# Restore the active-cell selection to match the saved workbook state.
$ws.Range("N18").Select()
